# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" column (E16:E46) currently lists the 31 contribution
# periods in descending order (2210 .. 2004). The update re-lists them in
# ascending order (2004 .. 2210), i.e. the chronological order used going
# forward for the new "estado de cuenta" periods being appended to the
# database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @(
    "2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
